$wb = $excel.ActiveWorkbook

# --- Sheet: VENTAS POR GRUPO ---
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws1.Range("L4").Value = 1216.13
$ws1.Range("M4").Value = 6518.91
$ws1.Range("H7").Value = 205.2
$ws1.Range("D12").Value = 475.2
$ws1.Range("C13").Value = 388.8
$ws1.Range("L13").Value = 886.88
$ws1.Range("M13").Value = 3133.95
$ws1.Range("H16").Value = 738
$ws1.Range("L16").Value = 1188
$ws1.Range("M16").Value = 6906.47
$ws1.Range("M17").Value = 1128.58
$ws1.Range("I22").Value = 27
$ws1.Range("M29").Value = 777.79
$ws1.Range("L32").Value = 2709.37
$ws1.Range("M32").Value = 4327.61
$ws1.Range("M37").Value = 2013.7
$ws1.Range("M52").Value = 1850.51

# Row 56 summary counters ("X de 54") - recompute as static labels
$ws1.Range("C56").Value = "3 de 54"
$ws1.Range("D56").Value = "7 de 54"
$ws1.Range("E56").Value = "2 de 54"
$ws1.Range("F56").Value = "0 de 54"
$ws1.Range("G56").Value = "0 de 54"
$ws1.Range("H56").Value = "3 de 54"
$ws1.Range("I56").Value = "4 de 54"
$ws1.Range("J56").Value = "0 de 54"
$ws1.Range("K56").Value = "2 de 54"
$ws1.Range("L56").Value = "10 de 54"
$ws1.Range("M56").Value = "15 de 54"

# --- Sheet: VENTA MENSUAL ---
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")
$ws2.Range("F4").Value = 8879.96
$ws2.Range("F7").Value = 70.41
$ws2.Range("F12").Value = 2181.51
$ws2.Range("F13").Value = 4409.63
$ws2.Range("F16").Value = 8832.47
$ws2.Range("F17").Value = 1128.58
$ws2.Range("F22").Value = 128.52
$ws2.Range("F29").Value = 4593.36
$ws2.Range("F32").Value = 7036.98
$ws2.Range("F37").Value = 2471.62
$ws2.Range("F53").Value = 1850.51
$ws2.Range("F54").Value = 1850.51
$ws2.Range("F60").Value = 56484.64

# --- Sheet: CUMPLIMIENTO MENSUAL ---
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$ws3.Range("D2").Value = 1296
$ws3.Range("E2").Value = 1530.66
$ws3.Range("F2").Value = 0.4584916473859608
$ws3.Range("D3").Value = 4294.18
$ws3.Range("E3").Value = 2329.08
$ws3.Range("F3").Value = 0.6483483964090192
$ws3.Range("D6").Value = 1725.3
$ws3.Range("E6").Value = 1182.28368146026
$ws3.Range("F6").Value = 0.5933793104566855
$ws3.Range("D7").Value = 442.8
$ws3.Range("E7").Value = 877.2
$ws3.Range("F7").Value = 0.3354545454545454
$ws3.Range("D11").Value = 15584.09
$ws3.Range("E11").Value = -1348.1
$ws3.Range("F11").Value = 1.09469661049214
$ws3.Range("D12").Value = 30178.52
$ws3.Range("E12").Value = 34765.48
$ws3.Range("F12").Value = 0.4646852673072185
$ws3.Range("D14").Value = 54610.73
$ws3.Range("E14").Value = 44345.52685923838
$ws3.Range("F14").Value = 0.5518673779029631
